# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el texto de conversión del día (celda A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nuevoTexto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.04 = 28486.19 pesos`n✅ 28486.19 pesos = 7.03 = 959.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $nuevoTexto

# --- tasas: actualizar valores de tasas en N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 141.999
$wsTasas.Range("O10").Value = 4045.01
$wsTasas.Range("N12").Value = 4053.63
$wsTasas.Range("O12").Value = 136.52
